$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.647003745318352
$ws1.Range("C2").Value = 0.5905420991926182
$ws1.Range("D2").Value = 0.9588014981273408
$ws1.Range("E2").Value = 0.7309064953604568
$ws1.Range("F2").Value = 0.8524808524808525
$ws1.Range("G2").Value = 0.9363438137441091
$ws1.Range("H2").Value = 0.7865799772755966
$ws1.Range("I2").Value = 512
$ws1.Range("J2").Value = 355
$ws1.Range("K2").Value = 179
$ws1.Range("L2").Value = 22

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.8905472636815921
$ws2.Range("C2").Value = 0.3352059925093633
$ws2.Range("D2").Value = 0.4870748299319728

$ws2.Range("B3").Value = 0.5905420991926182
$ws2.Range("C3").Value = 0.9588014981273408
$ws2.Range("D3").Value = 0.7309064953604568

$ws2.Range("B4").Value = 0.647003745318352
$ws2.Range("C4").Value = 0.647003745318352
$ws2.Range("D4").Value = 0.647003745318352
$ws2.Range("E4").Value = 0.647003745318352

$ws2.Range("B5").Value = 0.7405446814371052
$ws2.Range("C5").Value = 0.647003745318352
$ws2.Range("D5").Value = 0.6089906626462148

$ws2.Range("B6").Value = 0.7405446814371052
$ws2.Range("C6").Value = 0.647003745318352
$ws2.Range("D6").Value = 0.6089906626462148

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 179
$ws3.Range("C2").Value = 355

$ws3.Range("B3").Value = 22
$ws3.Range("C3").Value = 512
